# Updated cryptos list with latest price / volume(1h) figures from GitHub Actions.
# Prices in column D are stored as text (matching the source data feed),
# so a leading apostrophe forces plain-number-looking values to stay text
# instead of being auto-converted to floating point numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "51.736.62"
$ws.Range("E2").Value = "  +4.52%  "
$ws.Range("D3").Value = "2.770.27"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'117.06"
$ws.Range("E5").Value = "  +4.16%  "
$ws.Range("D6").Value = "'333.68"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +6.08%  "
$ws.Range("D10").Value = "'42.22"
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("D12").Value = "'20.27"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("E14").Value = "  +4.88%  "
$ws.Range("D15").Value = "3.202.79"
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").Value = "2.783.33"
$ws.Range("E16").Value = "  +5.89%  "
$ws.Range("E17").Value = "  +4.64%  "
$ws.Range("D18").Value = "51.672.53"
$ws.Range("E18").Value = "  +4.58%  "
$ws.Range("D19").Value = "'3.31"
$ws.Range("E19").Value = "  +12.65%  "
$ws.Range("D20").Value = "'13.59"
$ws.Range("E20").Value = "  +5.27%  "
$ws.Range("D21").Value = "'6.87"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D23").Value = "'279.25"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").Value = "'69.98"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "'2.71"
$ws.Range("E25").Value = "  +6.81%  "
$ws.Range("D26").Value = "'26.85"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'10.21"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +3.74%  "
$ws.Range("D31").Value = "'35.18"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "'50.02"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "'5.60"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "'0.0823"
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").Value = "'19.28"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").Value = "'5.04"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  +4.76%  "
$ws.Range("E40").Value = "  +9.05%  "
$ws.Range("D41").Value = "'128.29"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("D42").Value = "'23.27"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("D43").Value = "'2.32"
$ws.Range("E43").Value = "  +7.84%  "
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  +16.84%  "
$ws.Range("D46").Value = "2.089.97"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "'3.32"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("D48").Value = "'2.25"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("E49").Value = "  +6.72%  "
$ws.Range("D50").Value = "'60.87"
$ws.Range("E50").Value = "  +2.77%  "
$ws.Range("D51").Value = "'8.85"
$ws.Range("E51").Value = "  -0.74%  "
